$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1111
$ws.Range("E2").Value = 146
$ws.Range("F2").Value = 146
$ws.Range("G2").Value = 122
$ws.Range("H2").Value = 94
$ws.Range("I2").Value = 94
$ws.Range("K2").Value = 1638
$ws.Range("L2").Value = 869
$ws.Range("M2").Value = 769
$ws.Range("N2").Value = 769
$ws.Range("P2").Value = 396
$ws.Range("Q2").Value = 138
$ws.Range("R2").Value = -113
$ws.Range("S2").Value = -26
$ws.Range("T2").Value = 77
$ws.Range("U2").Value = 62
$ws.Range("V2").Value = 619
$ws.Range("W2").Value = 13.15
$ws.Range("X2").Value = 8.42
$ws.Range("Y2").Value = 12.78
$ws.Range("Z2").Value = 5.88
$ws.Range("AA2").Value = 112.95
$ws.Range("AB2").Value = 94.36
$ws.Range("AC2").Value = 591
$ws.Range("AE2").Value = 4859
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 386
$ws.Range("AI2").Value = 65.27
$ws.Range("AJ2").Value = 15834554
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3
$ws.Range("D3").Value = 1260
$ws.Range("E3").Value = 187
$ws.Range("F3").Value = 187
$ws.Range("G3").Value = 160
$ws.Range("H3").Value = 129
$ws.Range("I3").Value = 129
$ws.Range("K3").Value = 1799
$ws.Range("L3").Value = 970
$ws.Range("M3").Value = 828
$ws.Range("N3").Value = 828
$ws.Range("P3").Value = 396
$ws.Range("Q3").Value = 259
$ws.Range("R3").Value = -247
$ws.Range("S3").Value = -4
$ws.Range("T3").Value = 231
$ws.Range("U3").Value = 28
$ws.Range("V3").Value = 676
$ws.Range("W3").Value = 14.81
$ws.Range("X3").Value = 10.21
$ws.Range("Y3").Value = 16.1
$ws.Range("Z3").Value = 7.49
$ws.Range("AA3").Value = 117.16
$ws.Range("AB3").Value = 109.24
$ws.Range("AC3").Value = 812
$ws.Range("AE3").Value = 5231
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 404
$ws.Range("AI3").Value = 49.75
$ws.Range("AJ3").Value = 15834554
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()

# Row 4
$ws.Range("D4").Value = 1323
$ws.Range("E4").Value = 216
$ws.Range("F4").Value = 216
$ws.Range("G4").Value = 192
$ws.Range("H4").Value = 152
$ws.Range("I4").Value = 152
$ws.Range("K4").Value = 1922
$ws.Range("L4").Value = 1019
$ws.Range("M4").Value = 903
$ws.Range("N4").Value = 903
$ws.Range("P4").Value = 396
$ws.Range("Q4").Value = 161
$ws.Range("R4").Value = -113
$ws.Range("S4").Value = -53
$ws.Range("T4").Value = 134
$ws.Range("U4").Value = 27
$ws.Range("V4").Value = 687
$ws.Range("W4").Value = 16.35
$ws.Range("X4").Value = 11.46
$ws.Range("Y4").Value = 17.51
$ws.Range("Z4").Value = 8.15
$ws.Range("AA4").Value = 112.89
$ws.Range("AB4").Value = 128.11
$ws.Range("AC4").Value = 958
$ws.Range("AD4").Value = 14.28
$ws.Range("AE4").Value = 5703
$ws.Range("AF4").Value = 2.4
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 2.19
$ws.Range("AI4").Value = 31.33
$ws.Range("AJ4").Value = 15834554
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 1436
$ws.Range("E5").Value = 238
$ws.Range("F5").Value = 238
$ws.Range("G5").Value = 215
$ws.Range("H5").Value = 168
$ws.Range("I5").Value = 168
$ws.Range("K5").Value = 1985
$ws.Range("L5").Value = 964
$ws.Range("M5").Value = 1021
$ws.Range("N5").Value = 1021
$ws.Range("P5").Value = 396
$ws.Range("Q5").Value = 220
$ws.Range("R5").Value = -90
$ws.Range("S5").Value = -131
$ws.Range("T5").Value = 88
$ws.Range("U5").Value = 133
$ws.Range("V5").Value = 604
$ws.Range("W5").Value = 16.56
$ws.Range("X5").Value = 11.69
$ws.Range("Y5").Value = 17.45
$ws.Range("Z5").Value = 8.59
$ws.Range("AA5").Value = 94.45
$ws.Range("AB5").Value = 157.92
$ws.Range("AC5").Value = 1060
$ws.Range("AD5").Value = 19.86
$ws.Range("AE5").Value = 6448
$ws.Range("AF5").Value = 3.26
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 2.38
$ws.Range("AI5").Value = 47.18
$ws.Range("AJ5").Value = 15834554
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 1552
$ws.Range("E6").Value = 235
$ws.Range("F6").Value = 235
$ws.Range("G6").Value = 213
$ws.Range("H6").Value = 158
$ws.Range("I6").Value = 158
$ws.Range("K6").Value = 2009
$ws.Range("L6").Value = 916
$ws.Range("M6").Value = 1093
$ws.Range("N6").Value = 1093
$ws.Range("P6").Value = 396
$ws.Range("Q6").Value = 361
$ws.Range("R6").Value = -104
$ws.Range("S6").Value = -133
$ws.Range("T6").Value = 100
$ws.Range("U6").Value = 261
$ws.Range("V6").Value = 550
$ws.Range("W6").Value = 15.17
$ws.Range("X6").Value = 10.19
$ws.Range("Y6").Value = 14.97
$ws.Range("Z6").Value = 7.92
$ws.Range("AA6").Value = 83.84
$ws.Range("AB6").Value = 176
$ws.Range("AC6").Value = 999
$ws.Range("AD6").Value = 21.97
$ws.Range("AE6").Value = 6900
$ws.Range("AF6").Value = 3.18
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 2.28
$ws.Range("AI6").Value = 50.04
$ws.Range("AJ6").Value = 15834554

# Row 7
$ws.Range("D7").Value = 1677
$ws.Range("E7").Value = 297
$ws.Range("G7").Value = 280
$ws.Range("H7").Value = 220
$ws.Range("I7").Value = 221
$ws.Range("K7").Value = 2150
$ws.Range("L7").Value = 920
$ws.Range("M7").Value = 1230
$ws.Range("N7").Value = 1230
$ws.Range("P7").Value = 400
$ws.Range("Q7").Value = 210
$ws.Range("R7").Value = -80
$ws.Range("S7").Value = -160
$ws.Range("T7").Value = 80
$ws.Range("U7").Value = 130
$ws.Range("W7").Value = 17.71
$ws.Range("X7").Value = 13.12
$ws.Range("Y7").Value = 19.03
$ws.Range("Z7").Value = 10.58
$ws.Range("AA7").Value = 74.8
$ws.Range("AC7").Value = 1396
$ws.Range("AD7").Value = 13.69
$ws.Range("AE7").Value = 7768
$ws.Range("AF7").Value = 2.46
$ws.Range("AG7").Value = 750
$ws.Range("AH7").Value = 3.93
$ws.Range("AI7").Value = 53.74

# Row 8
$ws.Range("D8").Value = 1914
$ws.Range("E8").Value = 341
$ws.Range("G8").Value = 330
$ws.Range("H8").Value = 260
$ws.Range("I8").Value = 260
$ws.Range("K8").Value = 2370
$ws.Range("L8").Value = 990
$ws.Range("M8").Value = 1380
$ws.Range("N8").Value = 1380
$ws.Range("P8").Value = 400
$ws.Range("Q8").Value = 310
$ws.Range("R8").Value = -70
$ws.Range("S8").Value = -90
$ws.Range("T8").Value = 70
$ws.Range("U8").Value = 240
$ws.Range("W8").Value = 17.82
$ws.Range("X8").Value = 13.58
$ws.Range("Y8").Value = 19.92
$ws.Range("Z8").Value = 11.5
$ws.Range("AA8").Value = 71.73999999999999
$ws.Range("AC8").Value = 1642
$ws.Range("AD8").Value = 11.63
$ws.Range("AE8").Value = 8715
$ws.Range("AF8").Value = 2.19
$ws.Range("AG8").Value = 1000
$ws.Range("AH8").Value = 5.24
$ws.Range("AI8").Value = 60.9

# Row 9
$ws.Range("D9").Value = 2210
$ws.Range("E9").Value = 410
$ws.Range("G9").Value = 420
$ws.Range("H9").Value = 330
$ws.Range("I9").Value = 330
$ws.Range("K9").Value = 2520
$ws.Range("L9").Value = 970
$ws.Range("M9").Value = 1550
$ws.Range("N9").Value = 1550
$ws.Range("P9").Value = 400
$ws.Range("Q9").Value = 350
$ws.Range("R9").Value = -30
$ws.Range("S9").Value = -230
$ws.Range("T9").Value = 30
$ws.Range("U9").Value = 320
$ws.Range("W9").Value = 18.55
$ws.Range("X9").Value = 14.93
$ws.Range("Y9").Value = 22.53
$ws.Range("Z9").Value = 13.5
$ws.Range("AA9").Value = 62.58
$ws.Range("AC9").Value = 2084
$ws.Range("AD9").Value = 9.16
$ws.Range("AE9").Value = 9789
$ws.Range("AF9").Value = 1.95
$ws.Range("AG9").Value = 800
$ws.Range("AH9").Value = 4.19
$ws.Range("AI9").Value = 38.39
